$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038642423396848
$ws.Range("D2").Value = 1.044516097984918
$ws.Range("E2").Value = 1.046499221168531
$ws.Range("F2").Value = 1.055492008821157
$ws.Range("I2").Value = 1.036868983144936
$ws.Range("J2").Value = 1.043738752650489
$ws.Range("K2").Value = 1.047286928701409
$ws.Range("L2").Value = 1.049264485977202
$ws.Range("M2").Value = 1.058232318552414
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039616177363465
$ws.Range("D3").Value = 1.04525444585178
$ws.Range("E3").Value = 1.047349578577585
$ws.Range("F3").Value = 1.056392440981088
$ws.Range("I3").Value = 1.037043910859589
$ws.Range("J3").Value = 1.044357262356669
$ws.Range("K3").Value = 1.047836720561547
$ws.Range("L3").Value = 1.049926397693937
$ws.Range("M3").Value = 1.058945980818657
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.040246583607187
$ws.Range("D4").Value = 1.045732098966551
$ws.Range("E4").Value = 1.047900433007013
$ws.Range("F4").Value = 1.056975590338033
$ws.Range("I4").Value = 1.037155301556385
$ws.Range("J4").Value = 1.044757192132871
$ws.Range("K4").Value = 1.048191698784552
$ws.Range("L4").Value = 1.050354670623626
$ws.Range("M4").Value = 1.059407642136564
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.04051168265867
$ws.Range("D5").Value = 1.045932877220013
$ws.Range("E5").Value = 1.048132158427586
$ws.Range("F5").Value = 1.05722086690573
$ws.Range("I5").Value = 1.03720169874556
$ws.Range("J5").Value = 1.044925252833738
$ws.Range("K5").Value = 1.048340745256467
$ws.Range("L5").Value = 1.050534708735926
$ws.Range("M5").Value = 1.059601693209534
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.040556198402182
$ws.Range("D6").Value = 1.04596658714069
$ws.Range("E6").Value = 1.048171074680906
$ws.Range("F6").Value = 1.057262056950497
$ws.Range("I6").Value = 1.037209463715757
$ws.Range("J6").Value = 1.044953466844552
$ws.Range("K6").Value = 1.04836575985368
$ws.Range("L6").Value = 1.050564937436921
$ws.Range("M6").Value = 1.059634273365297
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040250125575135
$ws.Range("D7").Value = 1.045734781883747
$ws.Range("E7").Value = 1.047903528760565
$ws.Range("F7").Value = 1.056978867263522
$ws.Range("I7").Value = 1.037155923214151
$ws.Range("J7").Value = 1.044759438043598
$ws.Range("K7").Value = 1.048193691084524
$ws.Range("L7").Value = 1.050357076333805
$ws.Range("M7").Value = 1.059410235181189
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038971440833103
$ws.Range("D8").Value = 1.044765647473418
$ws.Range("E8").Value = 1.046786475023407
$ws.Range("F8").Value = 1.055796207737536
$ws.Range("I8").Value = 1.036928472896356
$ws.Range("J8").Value = 1.04394783981248
$ws.Range("K8").Value = 1.047472893007044
$ws.Range("L8").Value = 1.04948818733812
$ws.Range("M8").Value = 1.058473529398671
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036720742773146
$ws.Range("D9").Value = 1.043057149233229
$ws.Range("E9").Value = 1.044822862425213
$ws.Range("F9").Value = 1.053716174823376
$ws.Range("I9").Value = 1.03651392792973
$ws.Range("J9").Value = 1.042515542268365
$ws.Range("K9").Value = 1.04619687624649
$ws.Range("L9").Value = 1.047956924980987
$ws.Range("M9").Value = 1.056822022952031
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035222014608788
$ws.Range("D10").Value = 1.041917722671085
$ws.Range("E10").Value = 1.043517078096146
$ws.Range("F10").Value = 1.05233223361457
$ws.Range("I10").Value = 1.036228362666383
$ws.Range("J10").Value = 1.041559283414112
$ws.Range("K10").Value = 1.045342310310747
$ws.Range("L10").Value = 1.046936032291316
$ws.Range("M10").Value = 1.055720476266042
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.034573470977231
$ws.Range("D11").Value = 1.041424254068163
$ws.Range("E11").Value = 1.042952455624236
$ws.Range("F11").Value = 1.05173364033933
$ws.Range("I11").Value = 1.036102535716332
$ws.Range("J11").Value = 1.041144893814154
$ws.Range("K11").Value = 1.044971364806084
$ws.Range("L11").Value = 1.046493975189004
$ws.Range("M11").Value = 1.055243380175141
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034332636175688
$ws.Range("D12").Value = 1.041240945767738
$ws.Range("E12").Value = 1.042742849758554
$ws.Range("F12").Value = 1.051511396983049
$ws.Range("I12").Value = 1.036055471810863
$ws.Range("J12").Value = 1.04099092339964
$ws.Range("K12").Value = 1.044833443002753
$ws.Range("L12").Value = 1.046329775916221
$ws.Range("M12").Value = 1.055066148667759
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034384293187135
$ws.Range("D13").Value = 1.041280266546077
$ws.Range("E13").Value = 1.042787805475136
$ws.Range("F13").Value = 1.051559064353082
$ws.Range("I13").Value = 1.036065581928137
$ws.Range("J13").Value = 1.041023952718935
$ws.Range("K13").Value = 1.04486303384154
$ws.Range("L13").Value = 1.046364997182809
$ws.Range("M13").Value = 1.055104166168982
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.034553562204236
$ws.Range("D14").Value = 1.041409101985769
$ws.Range("E14").Value = 1.042935127081541
$ws.Range("F14").Value = 1.051715267569474
$ws.Range("I14").Value = 1.036098652049374
$ws.Range("J14").Value = 1.041132167538892
$ws.Range("K14").Value = 1.04495996691789
$ws.Range("L14").Value = 1.046480402411642
$ws.Range("M14").Value = 1.055228730499237
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.03465786282654
$ws.Range("D15").Value = 1.041488480190231
$ws.Range("E15").Value = 1.043025912732169
$ws.Range("F15").Value = 1.051811522927956
$ws.Range("I15").Value = 1.036118984440454
$ws.Range("J15").Value = 1.041198835971372
$ws.Range("K15").Value = 1.045019672577577
$ws.Range("L15").Value = 1.046551507466972
$ws.Range("M15").Value = 1.055305476504573
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035265065132704
$ws.Range("D16").Value = 1.041950470788316
$ws.Range("E16").Value = 1.043554566997882
$ws.Range("F16").Value = 1.052371974347639
$ws.Range("I16").Value = 1.036236667610283
$ws.Range("J16").Value = 1.041586778353697
$ws.Range("K16").Value = 1.04536690959687
$ws.Range("L16").Value = 1.046965370169449
$ws.Range("M16").Value = 1.055752137148386
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.035646058880332
$ws.Range("D17").Value = 1.04224024230194
$ws.Range("E17").Value = 1.043886390449609
$ws.Range("F17").Value = 1.052723708993159
$ws.Range("I17").Value = 1.036309905216252
$ws.Range("J17").Value = 1.041830038459856
$ws.Range("K17").Value = 1.045584478535729
$ws.Range("L17").Value = 1.047224974985123
$ws.Range("M17").Value = 1.056032284471853
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03586832613219
$ws.Range("D18").Value = 1.042409252409722
$ws.Range("E18").Value = 1.04408001366769
$ws.Range("F18").Value = 1.052928933684777
$ws.Range("I18").Value = 1.036352413544381
$ws.Range("J18").Value = 1.041971896719846
$ws.Range("K18").Value = 1.045711294675646
$ws.Range("L18").Value = 1.047376397645897
$ws.Range("M18").Value = 1.056195678099088
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.035944120290729
$ws.Range("D19").Value = 1.04246687897006
$ws.Range("E19").Value = 1.044146047070064
$ws.Range("F19").Value = 1.052998920786525
$ws.Range("I19").Value = 1.0363668721587
$ws.Range("J19").Value = 1.042020261393814
$ws.Range("K19").Value = 1.045754520720961
$ws.Range("L19").Value = 1.047428028755885
$ws.Range("M19").Value = 1.056251389107502
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.035605177709412
$ws.Range("D20").Value = 1.042209153452461
$ws.Range("E20").Value = 1.04385078104142
$ws.Range("F20").Value = 1.052685964583299
$ws.Range("I20").Value = 1.036302069223453
$ws.Range("J20").Value = 1.041803942177913
$ws.Range("K20").Value = 1.045561144557123
$ws.Range("L20").Value = 1.047197121883161
$ws.Range("M20").Value = 1.056002228503619
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.034503714947546
$ws.Range("D21").Value = 1.041371163482926
$ws.Range("E21").Value = 1.042891741216017
$ws.Range("F21").Value = 1.051669266821289
$ws.Range("I21").Value = 1.03608892272178
$ws.Range("J21").Value = 1.041100302280277
$ws.Range("K21").Value = 1.044931426292539
$ws.Range("L21").Value = 1.04644641842866
$ws.Range("M21").Value = 1.055192049862816
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.033811547500661
$ws.Range("D22").Value = 1.040844216641193
$ws.Range("E22").Value = 1.042289450064876
$ws.Range("F22").Value = 1.051030612656967
$ws.Range("I22").Value = 1.035953022230506
$ws.Range("J22").Value = 1.040657620605641
$ws.Range("K22").Value = 1.044534710338412
$ws.Range("L22").Value = 1.04597442427713
$ws.Range("M22").Value = 1.054682561882387
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03417844362686
$ws.Range("D23").Value = 1.041123567254314
$ws.Range("E23").Value = 1.042608669679067
$ws.Range("F23").Value = 1.05136911951019
$ws.Range("I23").Value = 1.036025244288312
$ws.Range("J23").Value = 1.040892320339576
$ws.Range("K23").Value = 1.044745091241384
$ws.Range("L23").Value = 1.046224636701035
$ws.Range("M23").Value = 1.054952659876683
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035623650027111
$ws.Range("D24").Value = 1.042223201192373
$ws.Range("E24").Value = 1.043866871164812
$ws.Range("F24").Value = 1.052703019460326
$ws.Range("I24").Value = 1.036305610619986
$ws.Range("J24").Value = 1.041815734061015
$ws.Range("K24").Value = 1.045571688450066
$ws.Range("L24").Value = 1.047209707502433
$ws.Range("M24").Value = 1.056015809537741
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03730229971511
$ws.Range("D25").Value = 1.043498918422168
$ws.Range("E25").Value = 1.045329929011062
$ws.Range("F25").Value = 1.05425343493183
$ws.Range("I25").Value = 1.036622722686599
$ws.Range("J25").Value = 1.042886075342691
$ws.Range("K25").Value = 1.04652744715261
$ws.Range("L25").Value = 1.048352806574322
$ws.Range("M25").Value = 1.057249077680434
